# Add a new "correction" set of columns: "Overtime 1.5" and "Overtime 2.0"
# These are inserted before the existing "Holiday OT" column, shifting the
# existing Holiday OT / Workday O / Restday OT columns two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at C:D (pushes old C,D,E -> E,F,G)
$ws.Columns("C:D").Insert()

# New header row 1 cells
$ws.Cells.Item(1, 3).Value = "Overtime 1.5"
$ws.Cells.Item(1, 4).Value = "Overtime 2.0"

# New data values for column C (Overtime 1.5)
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 3).Value = 76.59999999999999
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 3).Value = 63.1
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 3).Value = 57.72
$ws.Cells.Item(11, 3).Value = 0

# New data values for column D (Overtime 2.0)
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(8, 4).Value = 10
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(11, 4).Value = 0
